$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 4 through 6 entirely (data for the two "KHADIJA LALA" rows
# and the totals row are no longer needed), shifting remaining rows up.
$ws.Range("A4:K6").EntireRow.Delete() | Out-Null

# Row 2: replace "SAMIRA TATA" / "NABIL KAMAL"-style data with the
# "KHADIJA LALA" record (previously row 4 data).
$ws.Cells.Item(2, 1).Value = "KHADIJA LALA"
$ws.Cells.Item(2, 2).Value = "K5443645"
$ws.Cells.Item(2, 3).Value = "354564564324158786713544"
$ws.Cells.Item(2, 4).Value = "AG 100"
$ws.Cells.Item(2, 5).Value = "BP"
$ws.Cells.Item(2, 6).Value = "Direction régionale"
$ws.Cells.Item(2, 7).Value = "044/FES VILLE "
$ws.Cells.Item(2, 8).Value = "mensuelle"
$ws.Cells.Item(2, 9).Value = 20000
$ws.Cells.Item(2, 10).Value = 1500
$ws.Cells.Item(2, 11).Value = 18500

# Row 3: blank placeholder row (space-only text cells) carrying the same
# amounts as row 2, matching the previous row 5 amounts.
$ws.Cells.Item(3, 1).Value = " "
$ws.Cells.Item(3, 2).Value = " "
$ws.Cells.Item(3, 3).Value = " "
$ws.Cells.Item(3, 4).Value = " "
$ws.Cells.Item(3, 5).Value = " "
$ws.Cells.Item(3, 6).Value = " "
$ws.Cells.Item(3, 7).Value = " "
$ws.Cells.Item(3, 8).Value = " "
$ws.Cells.Item(3, 9).Value = 20000
$ws.Cells.Item(3, 10).Value = 1500
$ws.Cells.Item(3, 11).Value = 18500
